# Updates the "cryptos" sheet with refreshed price / volume(1h) figures.
# Columns: A=rank, B=Coin, C=Link, D=Price, E=Volume(1h)
#
# Price (column D) and Volume(1h) (column E) are stored as plain text in the
# workbook (e.g. "44.299.57" and "  +2.50%  "), not numbers. Several of the
# new Price values (e.g. "173.77") parse as valid numbers, so a naive
# `.Value = "..."` assignment would silently convert them into numeric
# cells. Set-TextValue forces the target range to Text format before the
# write and restores the "Normal" style afterwards so no stray number
# format / style index is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, [string]$value)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# row -> column -> new value. Only cells touched by the update are listed.
$updates = @{
    2  = @{ D = "44.417.31";  E = "  +2.73%  " }
    3  = @{ D = "2.368.15";   E = "  +1.01%  " }
    4  = @{               E = "  +0.10%  " }
    5  = @{ D = "0.678";      E = "  +4.41%  " }
    6  = @{ D = "239.14";     E = "  +3.35%  " }
    7  = @{ D = "73.63";      E = "  +11.33%  " }
    8  = @{               E = "  -0.03%  " }
    9  = @{ D = "0.551";      E = "  +20.95%  " }
    10 = @{               E = "  +7.62%  " }
    11 = @{ D = "29.57";      E = "  +10.47%  " }
    12 = @{ D = "0.108";      E = "  +2.82%  " }
    13 = @{ D = "2.716.53";   E = "  +0.86%  " }
    14 = @{ D = "16.88";      E = "  +10.47%  " }
    15 = @{ D = "6.74";       E = "  +7.74%  " }
    16 = @{ D = "0.909";      E = "  +9.02%  " }
    17 = @{ D = "2.334.54";   E = "  -0.62%  " }
    18 = @{ D = "44.276.59";  E = "  +2.43%  " }
    19 = @{ D = "0.0000103";  E = "  +5.67%  " }
    20 = @{ D = "78.18";      E = "  +6.27%  " }
    21 = @{ D = "6.47";       E = "  +4.75%  " }
    22 = @{ D = "256.15";     E = "  +3.61%  " }
    23 = @{               E = "  -0.02%  " }
    24 = @{               E = "  -3.82%  " }
    25 = @{               E = "  +3.62%  " }
    26 = @{               E = "  +6.74%  " }
    27 = @{               E = "  +1.07%  " }
    28 = @{ D = "22.54";      E = "  +1.17%  " }

    # Rows 29/30 swap places (Monero moves above ImmutableX) with refreshed values.
    29 = @{ B = "Monero";      C = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr";     D = "173.77"; E = "  -0.65%  " }
    30 = @{ B = "ImmutableX";  C = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx";      D = "1.59";   E = "  +5.73%  " }

    31 = @{               E = "  +3.93%  " }
    32 = @{               E = "  +5.73%  " }

    # Rows 33/34 swap places (Hedera moves above Filecoin) with refreshed values.
    33 = @{ B = "Hedera";      C = "https://coinranking.com/coin/jad286TjB+hedera-hbar";         D = "0.0743"; E = "  +8.22%  " }
    34 = @{ B = "Filecoin";    C = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil";         D = "5.23";   E = "  +5.56%  " }

    35 = @{ D = "5.23";       E = "  +4.99%  " }
    36 = @{               E = "  +10.13%  " }
    37 = @{               E = "  -1.52%  " }
    38 = @{ D = "6.50";       E = "  +0.79%  " }
    39 = @{               E = "  +7.86%  " }
    40 = @{               E = "  +10.52%  " }
    41 = @{               E = "  -0.05%  " }
    42 = @{               E = "  +0.09%  " }
    43 = @{               E = "  +4.13%  " }
    44 = @{ D = "0.0983";     E = "  +4.07%  " }
    45 = @{               E = "  +1.22%  " }
    46 = @{               E = "  +3.19%  " }
    47 = @{ D = "98.82";      E = "  +0.41%  " }
    48 = @{ D = "0.183";      E = "  +13.06%  " }
    49 = @{               E = "  +5.24%  " }
    50 = @{ D = "1.443.64";   E = "  +0.68%  " }
    51 = @{ D = "52.97";      E = "  +6.95%  " }
}

foreach ($row in $updates.Keys) {
    $rowData = $updates[$row]
    foreach ($col in $rowData.Keys) {
        $cell = $ws.Range("$col$row")
        if ($col -eq "D") {
            # Price column: some new values parse as numbers, force text.
            Set-TextValue $cell $rowData[$col]
        } else {
            $cell.Value = $rowData[$col]
        }
    }
}
